$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -18.49454702296344
$ws.Cells.Item(2, 3).Value = -18.49454702296344
$ws.Cells.Item(2, 4).Value = -18.49454702296344
$ws.Cells.Item(2, 5).Value = -18.49454702296344
$ws.Cells.Item(2, 6).Value = -18.49454702296344
$ws.Cells.Item(2, 7).Value = -18.49454702296344
$ws.Cells.Item(2, 8).Value = -18.49454702296344
$ws.Cells.Item(2, 9).Value = -18.49454702296344
$ws.Cells.Item(2, 10).Value = -18.49454702296344
$ws.Cells.Item(2, 11).Value = -18.49454702296344
$ws.Cells.Item(3, 2).Value = -18.49454702296344
$ws.Cells.Item(3, 3).Value = -18.49454702296344
$ws.Cells.Item(3, 4).Value = -18.49454702296344
$ws.Cells.Item(3, 5).Value = -18.49454702296344
$ws.Cells.Item(3, 6).Value = -18.49454702296344
$ws.Cells.Item(3, 7).Value = -18.49454702296344
$ws.Cells.Item(3, 8).Value = -18.49454702296344
$ws.Cells.Item(3, 9).Value = 1.245894209373187
$ws.Cells.Item(3, 10).Value = -18.49454702296344
$ws.Cells.Item(3, 11).Value = -18.49454702296344
$ws.Cells.Item(4, 2).Value = -18.49454702296344
$ws.Cells.Item(4, 3).Value = -18.49454702296344
$ws.Cells.Item(4, 4).Value = 1.565738586087327
$ws.Cells.Item(4, 5).Value = -18.49454702296344
$ws.Cells.Item(4, 6).Value = 3.509690248323069
$ws.Cells.Item(4, 7).Value = -18.49454702296344
$ws.Cells.Item(4, 8).Value = 1.473097195220993
$ws.Cells.Item(4, 9).Value = -18.49454702296344
$ws.Cells.Item(4, 10).Value = 0.926237680730864
$ws.Cells.Item(4, 11).Value = -18.49454702296344
$ws.Cells.Item(5, 2).Value = -18.49454702296344
$ws.Cells.Item(5, 3).Value = -18.49454702296344
$ws.Cells.Item(5, 4).Value = -18.49454702296344
$ws.Cells.Item(5, 5).Value = -18.49454702296344
$ws.Cells.Item(5, 6).Value = -18.49454702296344
$ws.Cells.Item(5, 7).Value = 2.84427693609596
$ws.Cells.Item(5, 8).Value = -18.49454702296344
$ws.Cells.Item(5, 9).Value = -18.49454702296344
$ws.Cells.Item(5, 10).Value = -18.49454702296344
$ws.Cells.Item(5, 11).Value = -18.49454702296344
$ws.Cells.Item(6, 2).Value = -18.49454702296344
$ws.Cells.Item(6, 3).Value = -18.49454702296344
$ws.Cells.Item(6, 4).Value = -18.49454702296344
$ws.Cells.Item(6, 5).Value = -18.49454702296344
$ws.Cells.Item(6, 6).Value = -18.49454702296344
$ws.Cells.Item(6, 7).Value = -18.49454702296344
$ws.Cells.Item(6, 8).Value = -18.49454702296344
$ws.Cells.Item(6, 9).Value = -18.49454702296344
$ws.Cells.Item(6, 10).Value = -18.49454702296344
$ws.Cells.Item(6, 11).Value = -18.49454702296344
$ws.Cells.Item(7, 2).Value = 2.386773847730157
$ws.Cells.Item(7, 3).Value = -18.49454702296344
$ws.Cells.Item(7, 4).Value = -18.49454702296344
$ws.Cells.Item(7, 5).Value = -18.49454702296344
$ws.Cells.Item(7, 6).Value = -18.49454702296344
$ws.Cells.Item(7, 7).Value = -18.49454702296344
$ws.Cells.Item(7, 8).Value = -18.49454702296344
$ws.Cells.Item(7, 9).Value = -18.49454702296344
$ws.Cells.Item(7, 10).Value = -18.49454702296344
$ws.Cells.Item(7, 11).Value = -18.49454702296344
$ws.Cells.Item(8, 2).Value = -18.49454702296344
$ws.Cells.Item(8, 3).Value = -18.49454702296344
$ws.Cells.Item(8, 4).Value = -18.49454702296344
$ws.Cells.Item(8, 5).Value = 1.62941900051691
$ws.Cells.Item(8, 6).Value = -18.49454702296344
$ws.Cells.Item(8, 7).Value = -18.49454702296344
$ws.Cells.Item(8, 8).Value = -18.49454702296344
$ws.Cells.Item(8, 9).Value = -18.49454702296344
$ws.Cells.Item(8, 10).Value = -18.49454702296344
$ws.Cells.Item(8, 11).Value = -18.49454702296344
$ws.Cells.Item(9, 2).Value = 3.884606296812005
$ws.Cells.Item(9, 3).Value = -18.49454702296344
$ws.Cells.Item(9, 4).Value = -18.49454702296344
$ws.Cells.Item(9, 5).Value = -18.49454702296344
$ws.Cells.Item(9, 6).Value = -18.49454702296344
$ws.Cells.Item(9, 7).Value = -18.49454702296344
$ws.Cells.Item(9, 8).Value = -18.49454702296344
$ws.Cells.Item(9, 9).Value = -18.49454702296344
$ws.Cells.Item(9, 10).Value = -18.49454702296344
$ws.Cells.Item(9, 11).Value = -18.49454702296344
$ws.Cells.Item(10, 2).Value = -18.49454702296344
$ws.Cells.Item(10, 3).Value = -18.49454702296344
$ws.Cells.Item(10, 4).Value = -18.49454702296344
$ws.Cells.Item(10, 5).Value = -18.49454702296344
$ws.Cells.Item(10, 6).Value = -18.49454702296344
$ws.Cells.Item(10, 7).Value = -18.49454702296344
$ws.Cells.Item(10, 8).Value = -18.49454702296344
$ws.Cells.Item(10, 9).Value = 1.748060739283616
$ws.Cells.Item(10, 10).Value = -18.49454702296344
$ws.Cells.Item(10, 11).Value = 2.211389886961461
$ws.Cells.Item(11, 2).Value = -18.49454702296344
$ws.Cells.Item(11, 3).Value = -18.49454702296344
$ws.Cells.Item(11, 4).Value = -18.49454702296344
$ws.Cells.Item(11, 5).Value = 3.069327457593031
$ws.Cells.Item(11, 6).Value = -18.49454702296344
$ws.Cells.Item(11, 7).Value = 2.854834103946436
$ws.Cells.Item(11, 8).Value = -18.49454702296344
$ws.Cells.Item(11, 9).Value = -18.49454702296344
$ws.Cells.Item(11, 10).Value = -18.49454702296344
$ws.Cells.Item(11, 11).Value = 1.943380767976893
$ws.Cells.Item(12, 2).Value = -18.49454702296344
$ws.Cells.Item(12, 3).Value = -18.49454702296344
$ws.Cells.Item(12, 4).Value = -18.49454702296344
$ws.Cells.Item(12, 5).Value = -18.49454702296344
$ws.Cells.Item(12, 6).Value = -18.49454702296344
$ws.Cells.Item(12, 7).Value = -18.49454702296344
$ws.Cells.Item(12, 8).Value = -18.49454702296344
$ws.Cells.Item(12, 9).Value = -18.49454702296344
$ws.Cells.Item(12, 10).Value = -18.49454702296344
$ws.Cells.Item(12, 11).Value = -18.49454702296344
$ws.Cells.Item(13, 2).Value = -18.49454702296344
$ws.Cells.Item(13, 3).Value = -18.49454702296344
$ws.Cells.Item(13, 4).Value = -18.49454702296344
$ws.Cells.Item(13, 5).Value = 2.537717353400387
$ws.Cells.Item(13, 6).Value = -18.49454702296344
$ws.Cells.Item(13, 7).Value = -18.49454702296344
$ws.Cells.Item(13, 8).Value = -18.49454702296344
$ws.Cells.Item(13, 9).Value = -18.49454702296344
$ws.Cells.Item(13, 10).Value = 1.681909822947921
$ws.Cells.Item(13, 11).Value = 1.777505198470952
$ws.Cells.Item(14, 2).Value = -18.49454702296344
$ws.Cells.Item(14, 3).Value = -18.49454702296344
$ws.Cells.Item(14, 4).Value = 1.616737680422731
$ws.Cells.Item(14, 5).Value = -18.49454702296344
$ws.Cells.Item(14, 6).Value = -18.49454702296344
$ws.Cells.Item(14, 7).Value = -18.49454702296344
$ws.Cells.Item(14, 8).Value = -18.49454702296344
$ws.Cells.Item(14, 9).Value = -18.49454702296344
$ws.Cells.Item(14, 10).Value = -18.49454702296344
$ws.Cells.Item(14, 11).Value = 1.957228537529889
$ws.Cells.Item(15, 2).Value = -18.49454702296344
$ws.Cells.Item(15, 3).Value = -18.49454702296344
$ws.Cells.Item(15, 4).Value = 2.003613643693139
$ws.Cells.Item(15, 5).Value = -18.49454702296344
$ws.Cells.Item(15, 6).Value = -18.49454702296344
$ws.Cells.Item(15, 7).Value = -18.49454702296344
$ws.Cells.Item(15, 8).Value = -18.49454702296344
$ws.Cells.Item(15, 9).Value = -18.49454702296344
$ws.Cells.Item(15, 10).Value = -18.49454702296344
$ws.Cells.Item(15, 11).Value = -18.49454702296344
$ws.Cells.Item(16, 2).Value = -18.49454702296344
$ws.Cells.Item(16, 3).Value = -18.49454702296344
$ws.Cells.Item(16, 4).Value = -18.49454702296344
$ws.Cells.Item(16, 5).Value = -18.49454702296344
$ws.Cells.Item(16, 6).Value = -18.49454702296344
$ws.Cells.Item(16, 7).Value = -18.49454702296344
$ws.Cells.Item(16, 8).Value = -18.49454702296344
$ws.Cells.Item(16, 9).Value = -18.49454702296344
$ws.Cells.Item(16, 10).Value = 1.877698463605407
$ws.Cells.Item(16, 11).Value = -18.49454702296344
$ws.Cells.Item(17, 2).Value = -18.49454702296344
$ws.Cells.Item(17, 3).Value = -18.49454702296344
$ws.Cells.Item(17, 4).Value = 2.026287433438187
$ws.Cells.Item(17, 5).Value = -18.49454702296344
$ws.Cells.Item(17, 6).Value = -18.49454702296344
$ws.Cells.Item(17, 7).Value = -18.49454702296344
$ws.Cells.Item(17, 8).Value = 2.07962239169333
$ws.Cells.Item(17, 9).Value = 2.108660037252444
$ws.Cells.Item(17, 10).Value = 2.56606777700902
$ws.Cells.Item(17, 11).Value = -18.49454702296344
$ws.Cells.Item(18, 2).Value = -18.49454702296344
$ws.Cells.Item(18, 3).Value = -18.49454702296344
$ws.Cells.Item(18, 4).Value = -18.49454702296344
$ws.Cells.Item(18, 5).Value = -18.49454702296344
$ws.Cells.Item(18, 6).Value = -18.49454702296344
$ws.Cells.Item(18, 7).Value = -18.49454702296344
$ws.Cells.Item(18, 8).Value = 1.964433242601521
$ws.Cells.Item(18, 9).Value = 2.045226149477592
$ws.Cells.Item(18, 10).Value = 2.404441773386488
$ws.Cells.Item(18, 11).Value = -18.49454702296344
$ws.Cells.Item(19, 2).Value = -18.49454702296344
$ws.Cells.Item(19, 3).Value = -18.49454702296344
$ws.Cells.Item(19, 4).Value = 2.039897405749873
$ws.Cells.Item(19, 5).Value = -18.49454702296344
$ws.Cells.Item(19, 6).Value = -18.49454702296344
$ws.Cells.Item(19, 7).Value = -18.49454702296344
$ws.Cells.Item(19, 8).Value = 1.673680696478157
$ws.Cells.Item(19, 9).Value = 1.800608687762933
$ws.Cells.Item(19, 10).Value = -18.49454702296344
$ws.Cells.Item(19, 11).Value = -18.49454702296344
$ws.Cells.Item(20, 2).Value = -18.49454702296344
$ws.Cells.Item(20, 3).Value = -18.49454702296344
$ws.Cells.Item(20, 4).Value = 0.8295390089827197
$ws.Cells.Item(20, 5).Value = -18.49454702296344
$ws.Cells.Item(20, 6).Value = 3.106012470887986
$ws.Cells.Item(20, 7).Value = -18.49454702296344
$ws.Cells.Item(20, 8).Value = 1.67189192401311
$ws.Cells.Item(20, 9).Value = 1.229761978177612
$ws.Cells.Item(20, 10).Value = -18.49454702296344
$ws.Cells.Item(20, 11).Value = 2.074233912771341
$ws.Cells.Item(21, 2).Value = -18.49454702296344
$ws.Cells.Item(21, 3).Value = 4.32192438393574
$ws.Cells.Item(21, 4).Value = -18.49454702296344
$ws.Cells.Item(21, 5).Value = 1.435926125059705
$ws.Cells.Item(21, 6).Value = -18.49454702296344
$ws.Cells.Item(21, 7).Value = 2.481374395307114
$ws.Cells.Item(21, 8).Value = 1.442351104432044
$ws.Cells.Item(21, 9).Value = -18.49454702296344
$ws.Cells.Item(21, 10).Value = -18.49454702296344
$ws.Cells.Item(21, 11).Value = -18.49454702296344
